$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.506.67'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '1.850.95'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''233.41'
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = '''0.4713'
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("D8").Value = '''0.2738'
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("D9").Value = '''0.06325'
$ws.Range("E9").Value = '  -1.70%  '
$ws.Range("D10").Value = '''17.62'
$ws.Range("E10").Value = '  +7.84%  '
$ws.Range("D11").Value = '1.844.60'
$ws.Range("E11").Value = '  -0.73%  '
$ws.Range("E12").Value = '  -0.80%  '
$ws.Range("D13").Value = '''5.042'
$ws.Range("E13").Value = '  +0.98%  '
$ws.Range("D14").Value = '''84.47'
$ws.Range("E14").Value = '  -1.28%  '
$ws.Range("D15").Value = '''0.6231'
$ws.Range("E15").Value = '  -1.54%  '
$ws.Range("D16").Value = '30.480.15'
$ws.Range("E16").Value = '  +0.57%  '
$ws.Range("D17").Value = '''241.73'
$ws.Range("E17").Value = '  +4.85%  '
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").Value = '''12.63'
$ws.Range("E19").Value = '  -0.95%  '
$ws.Range("D20").Value = '''0.000007333'
$ws.Range("E20").Value = '  -1.05%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").Value = '''4.923'
$ws.Range("E22").Value = '  -1.49%  '
$ws.Range("D23").Value = '''5.959'
$ws.Range("E23").Value = '  -0.59%  '
$ws.Range("D24").Value = '''9.201'
$ws.Range("E24").Value = '  -0.77%  '
$ws.Range("D25").Value = '''161.12'
$ws.Range("E25").Value = '  -3.14%  '
$ws.Range("D26").Value = '''17.98'
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").Value = '''1.878'
$ws.Range("E27").Value = '  -0.78%  '
$ws.Range("E28").Value = '  -2.73%  '
$ws.Range("D29").Value = '''1.362'
$ws.Range("E29").Value = '  -2.84%  '
$ws.Range("D30").Value = '''4.013'
$ws.Range("E30").Value = '  -3.34%  '
$ws.Range("D31").Value = '''3.829'
$ws.Range("E31").Value = '  -2.55%  '
$ws.Range("D32").Value = '''0.04856'
$ws.Range("E32").Value = '  -1.80%  '
$ws.Range("D33").Value = '''1.134'
$ws.Range("E33").Value = '  -2.69%  '
$ws.Range("D34").Value = '''0.7046'
$ws.Range("E34").Value = '  -2.82%  '
$ws.Range("D35").Value = '''2.713'
$ws.Range("E35").Value = '  +0.43%  '
$ws.Range("D36").Value = '''0.01902'
$ws.Range("E36").Value = '  +1.65%  '
$ws.Range("D37").Value = '''2.686'
$ws.Range("E37").Value = '  +1.37%  '
$ws.Range("D38").Value = '''0.8741'
$ws.Range("E38").Value = '  -4.49%  '
$ws.Range("D39").Value = '''1.971'
$ws.Range("E39").Value = '  -0.20%  '
$ws.Range("D40").Value = '''105.26'
$ws.Range("E40").Value = '  -0.83%  '
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").Value = '''0.4064'
$ws.Range("E42").Value = '  -1.19%  '
$ws.Range("D43").Value = '''5.495'
$ws.Range("E43").Value = '  -1.51%  '
$ws.Range("D44").Value = '''7.200'
$ws.Range("E44").Value = '  +1.03%  '
$ws.Range("D45").Value = '''62.22'
$ws.Range("E45").Value = '  +1.98%  '
$ws.Range("D46").Value = '''0.1211'
$ws.Range("E46").Value = '  +0.98%  '
$ws.Range("B47").Value = 'Elrond'
$ws.Range("C47").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D47").Value = '''33.31'
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''8.485'
$ws.Range("E48").Value = '  -2.39%  '
$ws.Range("D49").Value = '''0.05537'
$ws.Range("E49").Value = '  -0.86%  '
$ws.Range("D50").Value = '''1.365'
$ws.Range("E50").Value = '  -3.17%  '
$ws.Range("D51").Value = '''0.3666'
$ws.Range("E51").Value = '  -1.06%  '
